# Rebuild the "111_2" confirmations sheet with the new row layout:
#   - drop the old "Summary" section header row
#   - relabel each category's "New nominations"/"Carryover nominations"/
#     "Confirmed"/"Withdrawn"/"Returned to White House" rows with a
#     "<Category>, " prefix
#   - add two new summary rows ("Total new nominations" and
#     "Total carryover nominations") ahead of the existing totals
#   - the sheet shrinks from 44 rows to 43 rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the trailing row - the sheet goes from 44 rows to 43 rows overall.
$ws.Rows(44).Delete()

# --- Header / metadata rows (unchanged content, row positions keep their
#     existing formatting) ---
$ws.Range("A1").Value = "Labels"
$ws.Range("B1").Value = "Values"

$ws.Range("A2").Value = "Congress"
$ws.Range("B2").Value = 111

$ws.Range("A3").Value = "Session"
$ws.Range("B3").Value = 2

$ws.Range("A4").Value = "Start Date"
$ws.Range("B4").Value = 40183

$ws.Range("A5").Value = "End Date"
$ws.Range("B5").Value = 40543

# --- Civilian ---
$ws.Range("A6").Value = "Civilian "

$ws.Range("A7").Value = "     Civilian, New nominations"
$ws.Range("B7").Value = 435

$ws.Range("A8").Value = "     Civilian, Carryover nominations"
$ws.Range("B8").Value = 209

$ws.Range("A9").Value = "     Civilian, Confirmed "
$ws.Range("B9").Value = 453

$ws.Range("A10").Value = "     Civilian, Withdrawn "
$ws.Range("B10").Value = 16

$ws.Range("A11").Value = "     Civilian, Returned to White House "
$ws.Range("B11").Value = 175

# --- Other Civilian ---
$ws.Range("A12").Value = "Other Civilian "

$ws.Range("A13").Value = "     Other Civilian, New nominations"
$ws.Range("B13").Value = 2240

$ws.Range("A14").Value = "     Other Civilian, Carryover nominations"
$ws.Range("B14").Value = 112

$ws.Range("A15").Value = "     Other Civilian, Confirmed "
$ws.Range("B15").Value = 2347

$ws.Range("A16").Value = "     Other Civilian, Returned to White House "
$ws.Range("B16").Value = 5

# --- Air Force ---
$ws.Range("A17").Value = "Air Force "

$ws.Range("A18").Value = "     Air Force, New nominations"
$ws.Range("B18").Value = 6600

$ws.Range("A19").Value = "     Air Force, Carryover nominations"
$ws.Range("B19").Value = 759

$ws.Range("A20").Value = "     Air Force, Confirmed "
$ws.Range("B20").Value = 7318

$ws.Range("A21").Value = "     Air Force, Returned to White House "
$ws.Range("B21").Value = 41

# --- Army ---
$ws.Range("A22").Value = "Army "

$ws.Range("A23").Value = "     Army, New nominations"
$ws.Range("B23").Value = 7486

$ws.Range("A24").Value = "     Army, Carryover nominations"
$ws.Range("B24").Value = 76

$ws.Range("A25").Value = "     Army, Confirmed "
$ws.Range("B25").Value = 7553

$ws.Range("A26").Value = "     Army, Withdrawn "
$ws.Range("B26").Value = 5

$ws.Range("A27").Value = "     Army, Returned to White House "
$ws.Range("B27").Value = 4

# --- Navy ---
$ws.Range("A28").Value = "Navy "

$ws.Range("A29").Value = "     Navy, New nominations"
$ws.Range("B29").Value = 4448

$ws.Range("A30").Value = "     Navy, Carryover nominations"
$ws.Range("B30").Value = 8

$ws.Range("A31").Value = "     Navy, Confirmed "
$ws.Range("B31").Value = 4454

$ws.Range("A32").Value = "     Navy, Returned to White House "
$ws.Range("B32").Value = 2

# --- Marine Corps ---
$ws.Range("A33").Value = "Marine Corps "

$ws.Range("A34").Value = "     Marine Corps, New nominations"
$ws.Range("B34").Value = 627

$ws.Range("A35").Value = "     Marine Corps, Carryover nominations"
$ws.Range("B35").Value = 714

$ws.Range("A36").Value = "     Marine Corps, Confirmed "
$ws.Range("B36").Value = 1202

$ws.Range("A37").Value = "     Marine Corps, Returned to White House "
$ws.Range("B37").Value = 139

# --- Grand totals (replaces the old single "Summary" header row) ---
$ws.Range("A38").Value = "Total new nominations"
$ws.Range("B38").Value = 21836

$ws.Range("A39").Value = "Total carryover nominations"
$ws.Range("B39").Value = 1878

$ws.Range("A40").Value = "Total confirmed "
$ws.Range("B40").Value = 23327

$ws.Range("A41").Value = "Total unconfirmed "
$ws.Range("B41").Value = 0

$ws.Range("A42").Value = "Total withdrawn "
$ws.Range("B42").Value = 21

$ws.Range("A43").Value = "Total returned to the White House "
$ws.Range("B43").Value = 366

# Row 38 ("Total new nominations") needs the bold thousands-separator
# format already used by the other "Confirmed" subtotal rows (style with
# numFmtId 3, e.g. row 39's "Total carryover nominations"); row 41
# ("Total unconfirmed") needs to drop back to the plain/general format
# used elsewhere in column B (e.g. row 42's "Total withdrawn"). Copy the
# formatting from rows that already carry the desired style so Excel
# reuses the existing cellXf instead of minting a new one.
$ws.Range("B39").Copy()
$ws.Range("B38").PasteSpecial(-4122)

$ws.Range("B42").Copy()
$ws.Range("B41").PasteSpecial(-4122)

$excel.CutCopyMode = 0
